# Generate Report for Archive
#
# This localization-status report was refreshed:
#   - the "Ready for handoff" status (shown on the Overview sheet for
#     zh-cn/de-de, and on each language sheet's Status column) moved on
#     to "In Translation"
#   - a handoff archive run recorded a new handoff name,
#     "TestHandoff_2016-12-05-10-15", for the zh-cn and de-de rows that
#     previously had no "Lastest/Latest Handoff Name" value
#   - the Status and Lastest/Latest Handoff Name columns were re-sized to
#     fit the new text

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "Ready for handoff" -> "In Translation" -----------------------
# Overview sheet shows the per-language status in columns E (zh-cn) and F (de-de)
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# Each language sheet repeats the same status in its own "Status" column (C)
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# --- New handoff archive name ------------------------------------------------
# "Lastest/Latest Handoff Name" column (I) was blank, now records the new
# handoff archive name generated for this report
$wsZhCn.Range("I2").Value = "TestHandoff_2016-12-05-10-15"
$wsDeDe.Range("I2").Value = "TestHandoff_2016-12-05-10-15"

# --- Resize columns to fit the new content ----------------------------------
# Status column narrowed (shorter text than "Ready for handoff")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5

# Lastest/Latest Handoff Name column widened (now holds archive name text)
$wsZhCn.Columns.Item(9).ColumnWidth = 28.1666666666667
$wsDeDe.Columns.Item(9).ColumnWidth = 28.1666666666667
